$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Append the new row (row 58) with the 8 May 2020 data.
$ws.Range("A58").Value = 43959
$ws.Range("B58").Value = 33687
$ws.Range("C58").Value = 1848
$ws.Range("D58").Value = 48
$ws.Range("E58").Value = 3412

# Grow Table3 so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E58"))

$ws.Range("E57").Select()
